# Measure M13 ("Het project gebruikt ISO-25010 voor de specificatie van
# productkwaliteitseisen") has been deprecated - its guidance now lives
# under M01, so the whole slide introducing M13 is removed from the deck.

$p = $ppt.ActivePresentation

$targetIndex = -1
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    if ($slide.Shapes.Count -ge 1) {
        $title = $slide.Shapes.Item(1).TextFrame.TextRange.Text
        if ($title -like "M13:*") {
            $targetIndex = $i
            break
        }
    }
}

if ($targetIndex -ge 1) {
    $p.Slides.Item($targetIndex).Delete()
}
